$wb = $excel.ActiveWorkbook

# Add a new worksheet "GENBANK_GENOME" as the last sheet in the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "GENBANK_GENOME"

# Header row (row 1): column A is blank, B..L hold the ER sheet headers
$headers = @("TermSourceRef","Ontology","TAN","Content type (validation)","Notes during templating","Target term","Instruction","Requirement (m/o/n)","Value (cv/s/d)","Additional information","Review comments")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Data rows 2-16
$data = @(
    ,@("Source Name")
    ,@("Sample Name")
    ,@("Parameter [BioSample Accession Number]","NFDI4PSO:0000078","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000078")
    ,@("Parameter [Data filtering software]","NFDI4PSO:0000023","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000023")
    ,@("Parameter [Data filtering software version]","NFDI4PSO:0000024","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000024")
    ,@("Parameter [Data filtering Software Parameters]","NFDI4PSO:0000025","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000025")
    ,@("Parameter [Next generation sequencing instrument model]","NFDI4PSO:0000040","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000040")
    ,@("Parameter [sequence assembly algorithm]","OBI:0001522","OBI","http://purl.obolibrary.org/obo/OBI_0001522")
    ,@("Parameter [Sequence assembly algorithm version]","NFDI4PSO:0000060","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000060")
    ,@("Parameter [sequence assembly name]","OBI:0001948","OBI","http://purl.obolibrary.org/obo/OBI_0001948")
    ,@("Parameter [genome coverage]","OBI:0001939","OBI","http://purl.obolibrary.org/obo/OBI_0001939")
    ,@("Parameter [Genome status]","NFDI4PSO:0000061","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000061")
    ,@("Parameter [Genome reference sequence]","NFDI4PSO:0000026","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000026")
    ,@("Parameter [Processed data file name]","NFDI4PSO:0000028","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000028")
    ,@("Parameter [Processed data file format]","NFDI4PSO:0000027","NFDI4PSO","http://purl.obolibrary.org/obo/NFDI4PSO_0000027")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Column widths, matching the best-fit widths used in the authored sheet
$widths = @(54.7109375, 17.5703125, 9.85546875, 46.5703125, 23.5703125, 22.85546875, 11.28515625, 10.5703125, 20.42578125, 13.5703125, 21.5703125, 17.42578125)
for ($c = 0; $c -lt $widths.Length; $c++) {
    $ws.Columns.Item($c + 1).ColumnWidth = $widths[$c]
}

# Page margins: 2 cm top/bottom like the other sheets in this workbook
$ws.PageSetup.TopMargin = 56.692913385826778
$ws.PageSetup.BottomMargin = 56.692913385826778

# Select the entire sheet (all columns), matching the authored view state
$ws.Cells.Select()

Write-Host "done"
